# Auto commit at 2025-11-16  9:06:21.57
# Append two new daily rows (2025-11-15 data for both charging stations)
# to the bottom of the day-data log on Sheet1, then move the selection
# to just below the newly-added data (matching the author's saved cursor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: 四方坪站 (Sifangping station) for 2025-11-15 (serial 45976) ---
$ws.Range("A30").Value = 45976
$ws.Range("B30").Value = "四方坪站"
$ws.Range("C30").Value = 9342.76
$ws.Range("D30").Value = 8239.48
$ws.Range("E30").Value = 3054.69
$ws.Range("F30").Value = 384

# --- Row 31: 高岭站 (Gaoling station) for 2025-11-15 (serial 45976) ---
$ws.Range("A31").Value = 45976
$ws.Range("B31").Value = "高岭站"
$ws.Range("C31").Value = 4114.2
$ws.Range("D31").Value = 3535.09
$ws.Range("E31").Value = 1041.86
$ws.Range("F31").Value = 147

# Move / record the active selection the way the workbook was left
# (author's cursor ended up one row below the new data, in column H).
$ws.Range("H30").Select() | Out-Null
